# Events.xlsx — add "End of Battle Repair Attempt" gun-repair rows,
# inserted right after "e056c" (row 101) and before "e057" (old row 102).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# Minor copy tidy-up on the existing e022a "Rain Roll" entry (B36): extra
# spacing after the bold title to match the rest of the sheet's style.
$ws.Range("B36").Value = "<Bold>e022a Rain Roll</Bold>    
<InlineUIContainer><Button Content='Weather' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>"

# Insert 4 blank rows at 102..105 (everything from old row 102 onward shifts down by 4).
$ws.Rows("102:105").Insert()

$dMain = "<Bold>056d Repair Main Gun - End of Battle Repair Attempt</Bold> 
<InlineUIContainer><Button Content='r4.74.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
Attempt to repair malfunction gun by rolling on the <InlineUIContainer><Button Content='Gun Malfunction' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. Any result other than broken repairs the gun.
<LineBreak/><LineBreak/>"

$fCoax = "<Bold>056f Repair Coaxial MG - End of Battle Repair Attempt</Bold> 
<InlineUIContainer><Button Content='r4.74.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
Attempt to repair malfunction gun by rolling on the <InlineUIContainer><Button Content='Gun Malfunction' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. Any result other than broken repairs the gun.
<LineBreak/><LineBreak/>"

$eAA = "<Bold>056e Repair AA MG - End of Battle Repair Attempt</Bold> 
<InlineUIContainer><Button Content='r4.74.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
Attempt to repair malfunction gun by rolling on the <InlineUIContainer><Button Content='Gun Malfunction' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. Any result other than broken repairs the gun.
<LineBreak/><LineBreak/>"

$gBow = "<Bold>056g Repair Bow MG - End of Battle Repair Attempt</Bold> 
<InlineUIContainer><Button Content='r4.74.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<LineBreak/><LineBreak/>
Attempt to repair malfunction gun by rolling on the <InlineUIContainer><Button Content='Gun Malfunction' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Table. Any result other than broken repairs the gun.
<LineBreak/><LineBreak/>"

# Fill in the label/body cells in the same order the original author
# typed them (controls the order new entries land in the shared-string
# table): d's label+body together, then e/f/g labels, then f/e/g bodies.
$ws.Range("A102").Value = "e056d"
$ws.Range("B102").Value = $dMain
$ws.Range("A103").Value = "e056e"
$ws.Range("A104").Value = "e056f"
$ws.Range("A105").Value = "e056g"
$ws.Range("B104").Value = $fCoax
$ws.Range("B103").Value = $eAA
$ws.Range("B105").Value = $gBow

# All four new rows render the same as the other e056x rows: 90pt tall.
# (Cell styles - vertical-top on col A, wrap-text on col B - are already
# propagated from row 101 by Rows.Insert(), matching the rest of the table.)
$ws.Range("A102:A105").RowHeight = 90

# Scroll/selection state, matching where the author was working.
$ws.Range("B104").Select() | Out-Null
